$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand new data row at row 1318 (shifts existing rows 1318-1384 down to 1319-1385,
# and the worksheet dimension grows from A1:R1384 to A1:R1385).
$ws.Rows.Item(1318).Insert()

# Populate the newly inserted row with its data.
$ws.Range("A1318").Value = 9
$ws.Range("B1318").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C1318").Value = "Metropolitana"
$ws.Range("D1318").Value = 44516
$ws.Range("E1318").Value = 13
$ws.Range("F1318").Value = 100112002
$ws.Range("G1318").Value = "Pimiento"
$ws.Range("H1318").Value = "Cuatro cascos morado"
$ws.Range("I1318").Value = "Primera"
$ws.Range("J1318").Value = 25
$ws.Range("K1318").Value = 30000
$ws.Range("L1318").Value = 34000
$ws.Range("M1318").Value = 31920
$ws.Range("N1318").Value = "`$/caja 18 kilos"
$ws.Range("O1318").Value = "Limache"
$ws.Range("P1318").Value = 1773
$ws.Range("Q1318").Value = 18
$ws.Range("R1318").Value = "Hortaliza"
